$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / account holder info ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay TEXT (not become a number).
# Force text via NumberFormat, assign it, then restore the original cell
# format (style index) by pasting the format from a same-styled neighbour
# cell (B2), which avoids leaving a stray quotePrefix/numFmt behind.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 29.08.2025"

# --- Transaction row 6 ---
$ws.Range("B6").Value = "30.08."
$ws.Range("C6").Value = "31.08."
$ws.Range("D6").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E6").Value = "69,03-"

# --- Transaction row 7 ---
$ws.Range("B7").Value = "02.09."
$ws.Range("C7").Value = "03.09."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-15244846"
$ws.Range("E7").Value = "54,88-"

# --- Transaction row 8 ---
$ws.Range("B8").Value = "06.09."
$ws.Range("C8").Value = "07.09."
$ws.Range("D8").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E8").Value = "71,76-"

# --- Rows 9-11: no more transactions, cells become blank ---
$ws.Range("B9:D9").Value = ""
$ws.Range("B10:D10").Value = ""
$ws.Range("B11:D11").Value = ""

$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E9").VerticalAlignment = -4108    # xlCenter
$ws.Range("E9").WrapText = $true

$ws.Range("E10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152  # xlRight
$ws.Range("E10").VerticalAlignment = -4108    # xlCenter
$ws.Range("E10").WrapText = $true

$ws.Range("E11").Value = ""
$ws.Range("E11").HorizontalAlignment = -4152  # xlRight
$ws.Range("E11").VerticalAlignment = -4108    # xlCenter
$ws.Range("E11").WrapText = $true

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 11.09.2025"
$ws.Range("E12").Value = "195,67-"

# --- Next billing date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 21.09.2025"
